# Adds one-hot-encoding related bookkeeping to the feature preprocessing table:
#  - renames header J1 ("final_selection_correlation") to "selection_correlation_complete"
#    (done implicitly by writing the new text into J1)
#  - adds a new column K "selection_correlation_final" with "yes" markers mirroring
#    the existing selection columns for the rows that made the final cut
#  - marks the "dbsource" feature (row 16) as categorical instead of single_value
#  - adds a comment on row 7 (Arterial Blood Pressure mean) explaining weak correlation
#  - a handful of "yes" flags move between columns F/J/K as the selection got reviewed

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- header row -------------------------------------------------------
$ws.Range("J1").Value = "selection_correlation_complete"
$ws.Range("K1").Value = "selection_correlation_final"

# --- single cell value / content edits --------------------------------
$ws.Range("B16").Value = "categorical"
$ws.Range("H7").Value = "had little correlation, but should be influential"

$ws.Range("F7").Value = "yes"
$ws.Range("F29").Value = "yes"
$ws.Range("F43").Value = "yes"
$ws.Range("F60").Value = "yes"

$ws.Range("F37").ClearContents()
$ws.Range("F52").ClearContents()
$ws.Range("F63").ClearContents()

$ws.Range("J30").Value = "yes"
$ws.Range("J55").Value = "yes"

# --- new column K "yes" markers ---------------------------------------
$yesRows = @(3,5,9,12,13,14,15,16,17,18,19,20,21,23,25,26,27,28,31,32,33,34,35,39,40,41,44,47,48,49,50,53,58,64,65,66,67)
foreach ($r in $yesRows) {
    $ws.Range("K$r").Value = "yes"
}

# --- new column K cells that stay empty but carry the grey row banding
#     (mirrors the already-blank, grey-filled J cells on those rows) ---
$blankBandedRows = @(4,24,38,45)
foreach ($r in $blankBandedRows) {
    $srcRange = $ws.Range("J$r")
    $dstRange = $ws.Range("K$r")
    $srcRange.Copy()
    $dstRange.PasteSpecial(-4122)
    $dstRange.Value = ""
}

# --- copy the grey row-banding format into the new "yes" cells on the
#     banded rows, so column K keeps the same striping as columns A-J --
$bandedYesRows = @(12,15,16,17,18,19,20,21,23,25,26,34,39,40,41,44,53,58,65,66)
foreach ($r in $bandedYesRows) {
    $srcRange = $ws.Range("J$r")
    $dstRange = $ws.Range("K$r")
    $srcRange.Copy()
    $dstRange.PasteSpecial(-4122)
    $dstRange.Value = "yes"
}

# K1 header should look like the other bold header cells
$ws.Range("J1").Copy()
$ws.Range("K1").PasteSpecial(-4122)
$ws.Range("K1").Value = "selection_correlation_final"

$excel.CutCopyMode = 0

# --- restore the view: scrolled back to the top, editing near B17 -----
$ws.Range("B17").Select()
